# Weekly data refresh: insert the newest week's row for "Arándano (blue)"
# at the top of the historical series (row 216), pushing the rest of the
# rows down by one. This mirrors the author's "Fruta / hortaliza, semanal"
# commit: a new week of observations is prepended to the subset sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new weekly observation: insert a fresh row at 216.
$ws.Rows.Item(216).Insert()

# New week's record (2023-03-07 == serial 44992).
$ws.Range("A216").Value = 9
$ws.Range("B216").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C216").Value = "Metropolitana"
$ws.Range("D216").Value = 44992
$ws.Range("D216").NumberFormat = $ws.Range("D217").NumberFormat
$ws.Range("E216").Value = 13
$ws.Range("F216").Value = "Fruta"
$ws.Range("G216").Value = 100101
$ws.Range("H216").Value = "Berries"
$ws.Range("I216").Value = 100101001
$ws.Range("J216").Value = "Arándano (blue)"
$ws.Range("K216").Value = "Sin especificar"
$ws.Range("L216").Value = "Primera"
$ws.Range("M216").Value = 200
$ws.Range("N216").Value = 3000
$ws.Range("O216").Value = 3000
$ws.Range("P216").Value = 3000
$ws.Range("Q216").Value = "$/bandeja 2 kilos"
$ws.Range("R216").Value = "Provincia de Curicó"
$ws.Range("S216").Value = 1500
$ws.Range("T216").Value = 2
